$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.126.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.717.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.99%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3693'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.14'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3343'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.189'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07468'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.26%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.319'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.11'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.944'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.719.45'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001078'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06646'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.98'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.093'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.03'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '26.080.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.475'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.370'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +13.73%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.912.84'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.55'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.123'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.971'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08561'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.718'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.93'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.361'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06247'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02314'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.610'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2142'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.233'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.53'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +13.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6176'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.841'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5908'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.57'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.024'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07261'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.06'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.89%  '
